$d = $word.ActiveDocument

# 1. Title: merge the split runs "B" + "rief" + " Description on Case Study One"
#    into a single run with the same text "Brief Description on Case Study One".
$d.Content.Find.Execute("Brief Description on Case Study One", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Brief Description on Case Study One", 2) | Out-Null

# 2. Insert "different " before "class weights"
$d.Content.Find.Execute("many options with class weights", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "many options with different class weights", 2) | Out-Null

# 3. "is hard for me to have improve" -> "was hard for me to improve"
$d.Content.Find.Execute("It is hard for me to have improve", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "It was hard for me to improve", 2) | Out-Null

# 4. "I would attempt to add" -> "I would have attempted to add"
$d.Content.Find.Execute("I would attempt to add", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "I would have attempted to add", 2) | Out-Null
